$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$text = @'
questions = [
    {
        "title": "According to your company\u2019s accounting books, the closing balance of accounts receivable is $250,000 and the total credit sales for the year is $750,000. What is the accounts receivable turnover ratio?",
        "ques_type": 2,
        "options": [
            "3 times",
            "0.33 times",
            "1.5 times",
            "5 times"
        ],
        "score": "3 times"
    },
    {
        "title": "Cuckoos Company sold goods worth $350,000 to Grebes Company on 20 days\u2019 credit. Which of the following will be recorded in the books of Cuckoos Company for this transaction?",
        "ques_type": 2,
        "options": [
            "Debit Grebes Company: $350,000Credit Cuckoos Company: $350,000",
            "Debit Grebes Company: $350,000Credit sales: $350,000",
            "Debit sales: $350,000Credit Cuckoos Company: $350,000",
            "Debit cash: $350,000Credit sales: $350,000"
        ],
        "score": "Debit Grebes Company: $350,000Credit sales: $350,000"
    },
    {
        "title": "Your client Zoe purchased goods worth $50,000 during the year. Out of these purchases, 60% were on credit. During the year, Zoe paid $10,000 and also used a cash discount of $2,000. The opening balance was $4,000. What is Zoe\u2019s closing balance?",
        "ques_type": 2,
        "options": [
            "$40,000",
            "$54,000",
            "$22,000",
            "$46,000"
        ],
        "score": "$22,000"
    },
    {
        "title": "Which of the following balances are included in the list used to reconcile accounts receivable with the receivables ledger control account?",
        "ques_type": 2,
        "options": [
            "Opening balances of individual ledger accounts of accounts payable.",
            "Opening balances of individual ledger accounts of accounts receivable.",
            "Closing balances of individual ledger accounts of accounts receivable.",
            "Balances of all assets appearing on the balance sheet."
        ],
        "score": "Closing balances of individual ledger accounts of accounts receivable."
    }
]
'@

$ws.Range("A2").ClearContents()
$ws.Range("A1").ClearFormats()
$ws.Range("A1").Value = $text
